$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Pre outcome measures" column (column B) rating text.
# "A little stressful " is replaced by "Not stressful" for most rows,
# while row 6 is updated to the new "Moderately stressful" rating.
$ws.Range("B2").Value = "Not stressful"
$ws.Range("B3").Value = "Not stressful"
$ws.Range("B4").Value = "Not stressful"
$ws.Range("B5").Value = "Not stressful"
$ws.Range("B6").Value = "Moderately stressful"
$ws.Range("B7").Value = "Not stressful"
